# Hortaliza / Femacal de La Calera - Choclo
# Insert 4 new weekly records (rows 859-862) right before the existing
# row that used to be 859 ("Femacal de La Calera" / Coquimbo / Choclo),
# shifting the remaining records (old 859-885) down to 863-889.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows at 859; this shifts old rows 859:885 -> 863:889
# and also extends the sheet dimension to A1:R889 automatically.
$ws.Rows("859:862").Insert()

# Common / constant columns for this sheet (same market, region, category).
$commonA = 3
$commonB = "Femacal de La Calera"
$commonC = "Coquimbo"
$commonE = 5
$commonF = 100112024
$commonG = "Choclo"
$commonR = "Hortaliza"

# New row 859
$r = 859
$ws.Cells.Item($r, 1).Value = $commonA
$ws.Cells.Item($r, 2).Value = $commonB
$ws.Cells.Item($r, 3).Value = $commonC
$ws.Cells.Item($r, 4).Value = 44939
$ws.Cells.Item($r, 5).Value = $commonE
$ws.Cells.Item($r, 6).Value = $commonF
$ws.Cells.Item($r, 7).Value = $commonG
$ws.Cells.Item($r, 8).Value = "Choclero"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 27500
$ws.Cells.Item($r, 11).Value = 300
$ws.Cells.Item($r, 12).Value = 320
$ws.Cells.Item($r, 13).Value = 310
$ws.Cells.Item($r, 14).Value = "$/unidad"
$ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 16).Value = 310
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = $commonR

# New row 860
$r = 860
$ws.Cells.Item($r, 1).Value = $commonA
$ws.Cells.Item($r, 2).Value = $commonB
$ws.Cells.Item($r, 3).Value = $commonC
$ws.Cells.Item($r, 4).Value = 44939
$ws.Cells.Item($r, 5).Value = $commonE
$ws.Cells.Item($r, 6).Value = $commonF
$ws.Cells.Item($r, 7).Value = $commonG
$ws.Cells.Item($r, 8).Value = "Choclero"
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 3500
$ws.Cells.Item($r, 11).Value = 200
$ws.Cells.Item($r, 12).Value = 200
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = "$/unidad"
$ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 16).Value = 200
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = $commonR

# New row 861
$r = 861
$ws.Cells.Item($r, 1).Value = $commonA
$ws.Cells.Item($r, 2).Value = $commonB
$ws.Cells.Item($r, 3).Value = $commonC
$ws.Cells.Item($r, 4).Value = 44939
$ws.Cells.Item($r, 5).Value = $commonE
$ws.Cells.Item($r, 6).Value = $commonF
$ws.Cells.Item($r, 7).Value = $commonG
$ws.Cells.Item($r, 8).Value = "Dulce o Americano"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 23000
$ws.Cells.Item($r, 11).Value = 200
$ws.Cells.Item($r, 12).Value = 210
$ws.Cells.Item($r, 13).Value = 205
$ws.Cells.Item($r, 14).Value = "$/unidad"
$ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 16).Value = 205
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = $commonR

# New row 862
$r = 862
$ws.Cells.Item($r, 1).Value = $commonA
$ws.Cells.Item($r, 2).Value = $commonB
$ws.Cells.Item($r, 3).Value = $commonC
$ws.Cells.Item($r, 4).Value = 44939
$ws.Cells.Item($r, 5).Value = $commonE
$ws.Cells.Item($r, 6).Value = $commonF
$ws.Cells.Item($r, 7).Value = $commonG
$ws.Cells.Item($r, 8).Value = "Dulce o Americano"
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 2500
$ws.Cells.Item($r, 11).Value = 130
$ws.Cells.Item($r, 12).Value = 130
$ws.Cells.Item($r, 13).Value = 130
$ws.Cells.Item($r, 14).Value = "$/unidad"
$ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 16).Value = 130
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = $commonR
